$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- NOTA 2 text expanded ---
$ws.Range("B5").Value = "todas estas columnas SI o SI deben existir y NO cambiar el nombre de los títulos (case sensitive)"

# --- New row 8: NOTA 5 (label bold-italic-red, first so that font gets
#     registered before the plain-bold-red header font below, matching
#     the workbook's expected font creation order) ---
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "NOTA 5: "
$ws.Range("A8").Font.Color = 255

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "las columnas en rojo son opcionales, el resto es obligatorio (no pueden estar vacías)"
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").Font.Italic = $false

$ws.Range("B8").Characters(17, 4).Font.Bold = $true
$ws.Range("B8").Characters(17, 4).Font.Color = 255

# --- Header row: D1:F1 (name, lastname, email) turn bold + red ---
$ws.Range("D1:F1").Font.Color = 255

# --- I1/I2: custom date/time number format ---
$ws.Range("I1:I2").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"

# --- Column I width ---
$ws.Columns("I").ColumnWidth = 21.93

# --- Selection ---
[void]$ws.Range("J12").Select()
